# Carbon monoxide inhalation data update
# - Refresh several lab values in the data table (rows 4-9)
# - Update the active cell selection
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Arterial pO2) - updates to the 1 hr / 2 hr columns
$ws.Range("D4").Value = 23.5
$ws.Range("E4").Value = 42.3

# Row 5 (Arterial [O2])
$ws.Range("C5").Value = 83
$ws.Range("D5").Value = 77
$ws.Range("E5").Value = 70

# Row 6 (Venous pO2)
$ws.Range("C6").Value = 0.18
$ws.Range("E6").Value = 0.12

# Row 7 (Venous [O2])
$ws.Range("B7").Value = 41
$ws.Range("D7").Value = 27

# Row 9 (Cardiac Output)
$ws.Range("B9").Value = 5346
$ws.Range("C9").Value = 5468
$ws.Range("D9").Value = 5809
$ws.Range("E9").Value = 6493

# Update the selected cell in the sheet view
$ws.Range("E6").Select()
